$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF column (I) for rows 33-48 with the new value
$ws.Range("I33:I48").Value = 8.855
